$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 2883.5
$ws.Range("I40").Value = 2833.3333
$ws.Range("J40").Value = 2933.6667
$ws.Range("K40").Value = 2833.3333
$ws.Range("L40").Value = 2933.6667
$ws.Range("M40").Value = -2658.3333
$ws.Range("N40").Value = -3283.6667
# Row 53
$ws.Range("H53").Value = 294.07693
$ws.Range("I53").Value = 213.375
$ws.Range("J53").Value = 423.2
$ws.Range("K53").Value = 213.375
$ws.Range("L53").Value = 423.2
$ws.Range("M53").Value = 423.625
$ws.Range("N53").Value = -1697.2
# Row 130
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 29
$ws.Range("H29").Value = 8245
$ws.Range("I29").Value = 7500
$ws.Range("J29").Value = 8990
$ws.Range("K29").Value = 7500
$ws.Range("L29").Value = 8990
$ws.Range("M29").Value = -7192
$ws.Range("N29").Value = -9606
# Row 32
$ws.Range("H32").Value = 4385.66
$ws.Range("I32").Value = 2755.675
$ws.Range("J32").Value = 10905.6
$ws.Range("K32").Value = 2755.675
$ws.Range("L32").Value = 10905.6
$ws.Range("M32").Value = -2468.675
$ws.Range("N32").Value = -11479.6
# Row 45
$ws.Range("H45").Value = 1451.8334
$ws.Range("I45").Value = 1237
$ws.Range("J45").Value = 1666.6666
$ws.Range("K45").Value = 1237
$ws.Range("L45").Value = 1666.6666
$ws.Range("M45").Value = -860
$ws.Range("N45").Value = -2420.6666
# Row 61
$ws.Range("H61").Value = 3246.926
$ws.Range("I61").Value = 1935
$ws.Range("J61").Value = 5155.1816
$ws.Range("K61").Value = 1935
$ws.Range("L61").Value = 5155.1816
$ws.Range("M61").Value = -1723
$ws.Range("N61").Value = -5579.1816
# Row 127
$ws.Range("H127").Value = 41111.11
$ws.Range("J127").Value = 41111.11
$ws.Range("L127").Value = 41111.11
$ws.Range("N127").Value = -51031.11
# Row 132
$ws.Range("H132").Value = 3376.8965
$ws.Range("I132").Value = 2706.2104
$ws.Range("J132").Value = 4651.2
$ws.Range("K132").Value = 8118.6312
$ws.Range("L132").Value = 13953.6
$ws.Range("M132").Value = -5588.6312
$ws.Range("N132").Value = -19013.6
# Row 136
$ws.Range("H136").Value = 3246.926
$ws.Range("I136").Value = 1935
$ws.Range("J136").Value = 5155.1816
$ws.Range("K136").Value = 5805
$ws.Range("L136").Value = 15465.5448
$ws.Range("M136").Value = -3255
$ws.Range("N136").Value = -20565.5448

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 19
$ws.Range("H19").Value = 47500
$ws.Range("J19").Value = 45000
$ws.Range("L19").Value = 45000
$ws.Range("N19").Value = -45346
# Row 134
$ws.Range("H134").Value = 3448.25
$ws.Range("I134").Value = 2396.762
$ws.Range("J134").Value = 6602.7144
$ws.Range("K134").Value = 7190.286
$ws.Range("L134").Value = 19808.1432
$ws.Range("M134").Value = -4655.286
$ws.Range("N134").Value = -24878.1432

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 2237.6765
$ws.Range("I58").Value = 1235.3334
$ws.Range("J58").Value = 3856.8462
$ws.Range("K58").Value = 1235.3334
$ws.Range("L58").Value = 3856.8462
$ws.Range("M58").Value = -1032.3334
$ws.Range("N58").Value = -4262.8462
# Row 110
$ws.Range("H110").Value = 51900
$ws.Range("I110").Value = 63800
$ws.Range("J110").Value = 40000
$ws.Range("K110").Value = 63800
$ws.Range("L110").Value = 40000
$ws.Range("M110").Value = -59710
$ws.Range("N110").Value = -48180
# Row 134
$ws.Range("H134").Value = 3133.05
$ws.Range("I134").Value = 1413.4166
$ws.Range("J134").Value = 5712.5
$ws.Range("K134").Value = 4240.2498
$ws.Range("L134").Value = 17137.5
$ws.Range("M134").Value = -1705.2498
$ws.Range("N134").Value = -22207.5
# Row 136
$ws.Range("H136").Value = 2237.6765
$ws.Range("I136").Value = 1235.3334
$ws.Range("J136").Value = 3856.8462
$ws.Range("K136").Value = 3706.0002
$ws.Range("L136").Value = 11570.5386
$ws.Range("M136").Value = -1156.0002
$ws.Range("N136").Value = -16670.5386

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 26632.5
$ws.Range("I3").Value = 26632.5
$ws.Range("K3").Value = 79897.5
$ws.Range("M3").Value = -79785.5
# Row 113
$ws.Range("H113").Value = 16667527
$ws.Range("I113").Value = 556.25
$ws.Range("J113").Value = 19231676
$ws.Range("K113").Value = 1668.75
$ws.Range("L113").Value = 57695028
$ws.Range("M113").Value = 501.25
$ws.Range("N113").Value = -57699368
# Row 116
$ws.Range("H116").Value = 1017.8
$ws.Range("I116").Value = 772.25
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 2316.75
$ws.Range("L116").Value = 6000
$ws.Range("M116").Value = 1125.25
$ws.Range("N116").Value = -12884
# Row 133
$ws.Range("H133").Value = 12000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 12000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 36000
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -46120
# Row 134
$ws.Range("H134").Value = 3017.7334
$ws.Range("I134").Value = 2094.8462
$ws.Range("J134").Value = 9016.5
$ws.Range("K134").Value = 6284.5386
$ws.Range("L134").Value = 27049.5
$ws.Range("M134").Value = -1214.5386
$ws.Range("N134").Value = -37189.5
# Row 137
$ws.Range("H137").Value = 4595092.5
$ws.Range("I137").Value = 11114330
$ws.Range("J137").Value = 81774.46000000001
$ws.Range("K137").Value = 33342990
$ws.Range("L137").Value = 245323.38
$ws.Range("M137").Value = -33337890
$ws.Range("N137").Value = -255523.38
# Row 138
$ws.Range("H138").Value = 798.9286
$ws.Range("I138").Value = 798.9286
$ws.Range("K138").Value = 2396.7858
$ws.Range("M138").Value = 2743.2142

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2222.6667
$ws.Range("I80").Value = 2141.5386
$ws.Range("K80").Value = 2141.5386
$ws.Range("M80").Value = -1143.5386
# Row 83
$ws.Range("H83").Value = 2222.6667
$ws.Range("I83").Value = 2141.5386
$ws.Range("K83").Value = 10707.693
$ws.Range("M83").Value = -5715.692999999999
# Row 132
$ws.Range("H132").Value = 3049.3877
$ws.Range("I132").Value = 2398.8965
$ws.Range("J132").Value = 3992.6
$ws.Range("K132").Value = 7196.689499999999
$ws.Range("L132").Value = 11977.8
$ws.Range("M132").Value = -4666.689499999999
$ws.Range("N132").Value = -17037.8
# Row 134
$ws.Range("H134").Value = 33625.2
$ws.Range("J134").Value = 33625.2
$ws.Range("L134").Value = 100875.6
$ws.Range("N134").Value = -105945.6

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3245.7144
$ws.Range("I7").Value = 1980
$ws.Range("J7").Value = 3378.9473
$ws.Range("K7").Value = 1980
$ws.Range("L7").Value = 3378.9473
$ws.Range("M7").Value = -1868
$ws.Range("N7").Value = -3602.9473
# Row 16
$ws.Range("H16").Value = 499.46155
$ws.Range("I16").Value = 499.46155
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 499.46155
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -329.46155
$ws.Range("N16").ClearContents()
# Row 93
$ws.Range("H93").Value = 961
$ws.Range("I93").Value = 951.25
$ws.Range("J93").Value = 1000
$ws.Range("K93").Value = 951.25
$ws.Range("L93").Value = 1000
$ws.Range("M93").Value = 296.75
$ws.Range("N93").Value = -3496
# Row 126
$ws.Range("H126").Value = 3245.7144
$ws.Range("I126").Value = 1980
$ws.Range("J126").Value = 3378.9473
$ws.Range("K126").Value = 5940
$ws.Range("L126").Value = 10136.8419
$ws.Range("M126").Value = -3470
$ws.Range("N126").Value = -15076.8419
# Row 132
$ws.Range("H132").Value = 4094.432
$ws.Range("I132").Value = 3221.7727
$ws.Range("J132").Value = 4967.091
$ws.Range("K132").Value = 9665.3181
$ws.Range("L132").Value = 14901.273
$ws.Range("M132").Value = -7135.3181
$ws.Range("N132").Value = -19961.273
# Row 136
$ws.Range("H136").Value = 3193.5217
$ws.Range("I136").Value = 1694.3636
$ws.Range("J136").Value = 6999.077
$ws.Range("K136").Value = 5083.0908
$ws.Range("L136").Value = 20997.231
$ws.Range("M136").Value = -2533.0908
$ws.Range("N136").Value = -26097.231

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 26319666
$ws.Range("I132").Value = 45458070
$ws.Range("K132").Value = 136374210
$ws.Range("M132").Value = -136371680

Write-Host "All changes applied"